{"js": "// Insert a new list item right after the \"Namespace implicitly have public\n// access\" bullet, matching its paragraph/list formatting, with the text\n// \"Using static directive imports the members of a single class.\"\nconst body = context.document.body;\n\nconst results = body.search(\"Namespace implicitly have public access\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the anchor paragraph text in the document.\");\n}\n\nconst anchorParagraph = results.items[0].paragraphs.getFirst();\n\n// Inserting \"After\" an existing paragraph duplicates that paragraph's\n// formatting (style, numbering/list, spacing, indentation), matching how\n// Word behaves when pressing Enter at the end of a list item.\nanchorParagraph.insertParagraph(\n  \"Using static directive imports the members of a single class.\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# Insert a new list item right after the \"Namespace implicitly have public\n# access\" bullet, matching its paragraph/list formatting, with the text\n# \"Using static directive imports the members of a single class.\"\n$d = $word.ActiveDocument\n\n# Find the anchor paragraph by its text.\n$found = $d.Content\n$ok = $found.Find.Execute(\"Namespace implicitly have public access\")\nif (-not $ok) {\n    throw \"Could not find the anchor paragraph text in the document.\"\n}\n\n# Resolve the paragraph index of the found range.\n$preceding = $d.Range(0, $found.Start)\n$targetIndex = $preceding.Paragraphs.Count + 1\n\n# Inserting a new paragraph right after the anchor duplicates that\n# paragraph's formatting (style, numbering/list, spacing, indentation),\n# matching how Word behaves when pressing Enter at the end of a list item.\n$d.Paragraphs.Item($targetIndex).Range.InsertParagraphAfter()\n\n# Re-fetch the freshly created paragraph and set its text.\n$newPara = $d.Paragraphs.Item($targetIndex + 1)\n$newPara.Range.Text = \"Using static directive imports the members of a single class.\"\n"}
